$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet so the values below can be written
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer string (A80)
$disclaimerCell = $ws.Range("A80")
$disclaimerText = $disclaimerCell.Value()
$disclaimerCell.Value = $disclaimerText.Replace("2021-05-12", "2021-05-13")

# Refresh Weight (D) and Percent Change (E) figures for each holding
$ws.Range("D2").Value = 0.06242849816859322
$ws.Range("E2").Value = 0.01791968722000492
$ws.Range("D3").Value = 0.03816095923029043
$ws.Range("E3").Value = 0.003023534711955156
$ws.Range("D4").Value = 0.03182965469407208
$ws.Range("E4").Value = 0.01686192468619252
$ws.Range("D5").Value = 0.02938316617218879
$ws.Range("E5").Value = 0.01409599782441195
$ws.Range("D6").Value = 0.02663872108810654
$ws.Range("E6").Value = 0.01308487671855474
$ws.Range("D7").Value = 0.02541690950466158
$ws.Range("E7").Value = 0.02572245157192765
$ws.Range("D8").Value = 0.1889378761233708
$ws.Range("E8").Value = 0.02431834929992638
$ws.Range("D9").Value = 0.02511584544479373
$ws.Range("E9").Value = 0.01046373365041631
$ws.Range("D10").Value = 0.02297401276579817
$ws.Range("E10").Value = 0.01608381289656191
$ws.Range("D11").Value = 0.02208050629321449
$ws.Range("E11").Value = 0.01255282127765356
$ws.Range("D12").Value = 0.02009703492399393
$ws.Range("E12").Value = 0.002755130728141708
$ws.Range("D13").Value = 0.02027525193315967
$ws.Range("E13").Value = 0.01602719766877136
$ws.Range("D14").Value = 0.0172841441769712
$ws.Range("E14").Value = 0.01724806201550377
$ws.Range("D15").Value = 0.01615914930661244
$ws.Range("E15").Value = -0.000239757844576971
$ws.Range("D16").Value = 0.01444009773903454
$ws.Range("E16").Value = 0.01135804677368868
$ws.Range("D17").Value = 0.0142639792997231
$ws.Range("E17").Value = 0.01691922907164933
$ws.Range("D18").Value = 0.01455176586955121
$ws.Range("E18").Value = 0.01483741246619985
$ws.Range("D19").Value = 0.01343104885312623
$ws.Range("E19").Value = 0.008957197157494434
$ws.Range("D20").Value = 0.0135690298210469
$ws.Range("E20").Value = -0.01232511658894075
$ws.Range("D21").Value = 0.01277314720380106
$ws.Range("E21").Value = 0.007507037847982456
$ws.Range("D22").Value = 0.01345441562131168
$ws.Range("E22").Value = 0.01426583159359773
$ws.Range("D23").Value = 0.01151372279142892
$ws.Range("E23").Value = 0.009533957714093377
$ws.Range("D24").Value = 0.01308765018215899
$ws.Range("E24").Value = 0.007869353923575906
$ws.Range("D25").Value = 0.01151388421987925
$ws.Range("E25").Value = 0.02656151419558372
$ws.Range("D26").Value = 0.008678393489810142
$ws.Range("E26").Value = 0.005078124999999822
$ws.Range("D27").Value = 0.00900238038962504
$ws.Range("E27").Value = 0.04852332024315476
$ws.Range("D28").Value = 0.01002995319020887
$ws.Range("E28").Value = 0.01586931155192528
$ws.Range("D29").Value = 0.01004133389595722
$ws.Range("E29").Value = -0.001571467614102318
$ws.Range("D30").Value = 0.009717185567691992
$ws.Range("E30").Value = 0.01578204169781539
$ws.Range("D31").Value = 0.008439398669094668
$ws.Range("E31").Value = 0.007273405445729297
$ws.Range("D32").Value = 0.01051109068642127
$ws.Range("E32").Value = 0.01359943482868253
$ws.Range("D33").Value = 0.009664761678446905
$ws.Range("E33").Value = 0.006848142441362892
$ws.Range("D34").Value = 0.009159773128698048
$ws.Range("E34").Value = 0.008697261287935021
$ws.Range("D35").Value = 0.009354375125572419
$ws.Range("E35").Value = 0.0265757798006816
$ws.Range("D36").Value = 0.008451990088220507
$ws.Range("E36").Value = 0.009454232917920136
$ws.Range("D37").Value = 0.008698370260538638
$ws.Range("E37").Value = 0.01468439000626343
$ws.Range("D38").Value = 0.007141877142444342
$ws.Range("E38").Value = -0.03085320992049356
$ws.Range("D39").Value = 0.008969933271108451
$ws.Range("E39").Value = 0.01033005794910569
$ws.Range("D40").Value = 0.008232972038233531
$ws.Range("E40").Value = -0.006426376082704688
$ws.Range("D41").Value = 0.006818414885093021
$ws.Range("E41").Value = 0.01946114872863314
$ws.Range("D42").Value = 0.007041024718099869
$ws.Range("E42").Value = 0.01998074145402007
$ws.Range("D43").Value = 0.008053261815902222
$ws.Range("E43").Value = 0.0181107491856678
$ws.Range("D44").Value = 0.007536529346391763
$ws.Range("E44").Value = 0.001649299047904895
$ws.Range("D45").Value = 0.007341443064166399
$ws.Range("E45").Value = 0.004705571924886964
$ws.Range("D46").Value = 0.007930656907875607
$ws.Range("E46").Value = 0.02540302882266743
$ws.Range("D47").Value = 0.007554851475504365
$ws.Range("E47").Value = 0.008846153846153726
$ws.Range("D48").Value = 0.007145388211239049
$ws.Range("E48").Value = 0.02156404260846956
$ws.Range("D49").Value = 0.006646211085702117
$ws.Range("E49").Value = 0.006375808361417157
$ws.Range("D50").Value = 0.00740294730374262
$ws.Range("E50").Value = 0.01067402254737337
$ws.Range("D51").Value = 0.006704284970708799
$ws.Range("E51").Value = 0.01158772964773291
$ws.Range("D52").Value = 0.006646816442390859
$ws.Range("E52").Value = 0.01880995749848213
$ws.Range("D53").Value = 0.005209699663316552
$ws.Range("E53").Value = 0.02184522426214253
$ws.Range("D54").Value = 0.006237595320801039
$ws.Range("E54").Value = 0.009006211180124346
$ws.Range("D55").Value = 0.005590469020535481
$ws.Range("E55").Value = 0.02837033026529512
$ws.Range("D56").Value = 0.005705432309294539
$ws.Range("E56").Value = 0.01405318459274407
$ws.Range("D57").Value = 0.006802594896960554
$ws.Range("E57").Value = 0.02358803986710956
$ws.Range("D58").Value = 0.005540224415369869
$ws.Range("E58").Value = 0.01538461538461533
$ws.Range("D59").Value = 0.005392275240641241
$ws.Range("E59").Value = 0.01427994072477423
$ws.Range("D60").Value = 0.005030110512322992
$ws.Range("E60").Value = 0.0177150192554556
$ws.Range("D61").Value = 0.004885874191951984
$ws.Range("E61").Value = 0.004361257495911319
$ws.Range("D62").Value = 0.005047867641859433
$ws.Range("E62").Value = 0.02006715701950768
$ws.Range("D63").Value = 0.004293351065010984
$ws.Range("E63").Value = 0.00579034441269366
$ws.Range("D64").Value = 0.004166791159951253
$ws.Range("E64").Value = 0.007903300790330192
$ws.Range("D65").Value = 0.003900757073805288
$ws.Range("E65").Value = 0.01009766594934614
$ws.Range("D66").Value = 0.003833602838467471
$ws.Range("E66").Value = 0.01768569984840829
$ws.Range("D67").Value = 0.003789048586176035
$ws.Range("E67").Value = 0.03610685071574649
$ws.Range("D68").Value = 0.003583832668692383
$ws.Range("E68").Value = 0.02496537279145961
$ws.Range("D69").Value = 0.003662407966891138
$ws.Range("E69").Value = 0.00484848484848488
$ws.Range("D70").Value = 0.002931702086466555
$ws.Range("E70").Value = 0.01150817686250738
$ws.Range("D71").Value = 0.002904622463923481
$ws.Range("E71").Value = 0.003709724480013499
$ws.Range("D72").Value = 0.002221013333883051
$ws.Range("E72").Value = -0.006777628375186251
$ws.Range("D73").Value = 0.001957239246041724
$ws.Range("E73").Value = 0.003464060373623568
$ws.Range("D74").Value = 0.001905218927922465
$ws.Range("E74").Value = 0.004384757143764961
$ws.Range("D75").Value = 0.001400068949723277
$ws.Range("E75").Value = 0.01372074253430178
$ws.Range("D76").Value = 0.00171211014421366
$ws.Range("E76").Value = -0.00900433716764093
$ws.Range("E77").Value = 0.01476350954999539

# Restore sheet protection
$ws.Protect("D382")
